# Update BOM workbook to v1.2:
#  - U1 (MSP430FR5994) part numbers updated to the -IPNR/-IPN variants
#  - U2 (INA226-Q1) manufacturer part number set to INA226AQDGSRQ1
#  - U3 (MS5611-01BA01) manufacturer part number set to MS561101BA03-50
#  - U5 (TMP175C-Q1) manufacturer part number / description updated to TMP175AIDR / TMP75C-Q1 SOIC
#  - Iterative-calculation max change set to 1E-4
#  - Last active selection moved to C30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Iterative calculation "Maximum Change" -> calcPr iterateDelta="1E-4"
$excel.MaxChange = 0.0001

# Leading "'" forces these to stay text cells with the original quote-prefixed
# style (matches the unchanged s="3" style in the target file) instead of
# Excel reinterpreting/restyling the cell when the value is assigned.
$ws.Range("D24").Value = "'MSP430FR5994IPNR"
$ws.Range("E24").Value = "'MSP430FR5994IPN"

$ws.Range("D25").Value = "'INA226AQDGSRQ1"

$ws.Range("D26").Value = "'MS561101BA03-50"

$ws.Range("D28").Value = "'TMP175AIDR"
$ws.Range("E28").Value = "'TMP75C-Q1 SOIC"

# Restore the last selected cell recorded in the workbook view
$ws.Range("C30").Select() | Out-Null
